$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 693 is a new "section header" row (like row 692): column A gets the
# orange header fill, B/C stay unfilled. Copy that formatting first, then
# set the three shared-string values for the row.
$ws.Range("A692:C692").Copy()
$ws.Range("A693:C693").PasteSpecial(-4122)

# Rows 694-705 are "detail" rows (like row 649): same visual formatting as
# the header row (A = orange fill, B/C = no fill), just copy it down too.
$ws.Range("A649:C649").Copy()
$ws.Range("A694:C705").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("A693").Value = "referral.hoverDetails"
$ws.Range("B693").Value = "Hover over the details to see more information about each referral."
$ws.Range("C693").Value = "&Hover over the details to see more information about each referral."
$ws.Range("A694").Value = "statistics.totalFollowedUpFAdult"
$ws.Range("B694").Value = "Total Followed Up Female Adult Clients:"
$ws.Range("C694").Value = "&Total Followed Up Female Adult Clients:"
$ws.Range("A695").Value = "statistics.totalFollowedUpMAdult"
$ws.Range("B695").Value = "Total Followed Up Male Adult Clients:"
$ws.Range("C695").Value = "&Total Followed Up Male Adult Clients:"
$ws.Range("A696").Value = "statistics.totalFollowedUpFChild"
$ws.Range("B696").Value = "Total Followed Up Female Child Clients:"
$ws.Range("C696").Value = "&Total Followed Up Female Child Clients:"
$ws.Range("A697").Value = "statistics.totalFollowedUpMChild"
$ws.Range("B697").Value = "Total Followed Up Male Child Clients:"
$ws.Range("C697").Value = "&Total Followed Up Male Child Clients:"
$ws.Range("A698").Value = "statistics.followUpVisits"
$ws.Range("B698").Value = "Follow Up Visits:"
$ws.Range("C698").Value = "&Follow Up Visits:"
$ws.Range("A699").Value = "statistics.newClients"
$ws.Range("B699").Value = "New Clients:"
$ws.Range("C699").Value = "&New Clients:"
$ws.Range("A700").Value = "statistics.totalNewFAdult"
$ws.Range("B700").Value = "Total New Female Adult Clients:"
$ws.Range("C700").Value = "&Total New Female Adult Clients:"
$ws.Range("A701").Value = "statistics.totalNewMAdult"
$ws.Range("B701").Value = "Total New Male Adult Clients:"
$ws.Range("C701").Value = "&Total New Male Adult Clients:"
$ws.Range("A702").Value = "statistics.totalNewFChild"
$ws.Range("B702").Value = "Total New Female Child Clients:"
$ws.Range("C702").Value = "&Total New Female Child Clients:"
$ws.Range("A703").Value = "statistics.totalNewMChild"
$ws.Range("B703").Value = "Total New Male Child Clients:"
$ws.Range("C703").Value = "&Total New Male Child Clients:"
$ws.Range("A704").Value = "statistics.allChildren"
$ws.Range("B704").Value = "All Children "
$ws.Range("C704").Value = "&All Children "
$ws.Range("A705").Value = "statistics.allAdults"
$ws.Range("B705").Value = "All Adults"
$ws.Range("C705").Value = "&All Adults"
